$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Name) to hold the new "Sex" column
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B2").Value = "Sex"

# Data
$ws.Range("B3").Value = "Male"
$ws.Range("B4").Value = "Female"
$ws.Range("B5").Value = "Male"
